$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.868023157119751
$ws.Range("B1").Value = 2.4712073802948
$ws.Range("C1").Value = 1.901860117912292
$ws.Range("D1").Value = 1.832730174064636
$ws.Range("E1").Value = 1.783858418464661
